$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    foreach ($ftr in $sec.Footers) {
        $count = $ftr.Shapes.Count
        for ($i = $count; $i -ge 1; $i--) {
            $shp = $ftr.Shapes.Item($i)
            $shp.Delete()
        }
    }
}
